# "Take failed screenshot implemented"
# Trim the homepage test data providers down to two data rows each,
# and swap sheet2's password value to "admin".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("homePageTest1")
$ws2 = $wb.Worksheets.Item("homePageTest2")

# --- homePageTest1 (sheet1): drop the 4th (last) data row ---
$ws1.Rows.Item(4).Delete()
$ws1.Range("A4:XFD4").Select()

# --- homePageTest2 (sheet2): change the password for row 2, drop row 3 ---
$ws2.Cells.Item(2, 2).Value = "admin"
$ws2.Rows.Item(3).Delete()
$ws2.Range("H9").Select()
